$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 3428.5
$ws.Range("J32").Value = 2492.5
$ws.Range("L32").Value = 2492.5
$ws.Range("N32").Value = -3144.5
$ws.Range("H82").Value = 22229460
$ws.Range("J82").Value = 40011220
$ws.Range("L82").Value = 120033660
$ws.Range("N82").Value = -120034472
$ws.Range("H85").Value = 22229460
$ws.Range("J85").Value = 40011220
$ws.Range("L85").Value = 120033660
$ws.Range("N85").Value = -120036468
$ws.Range("H92").Value = 1809.7778
$ws.Range("I92").Value = 1905.4286
$ws.Range("K92").Value = 1905.4286
$ws.Range("M92").Value = -657.4286
$ws.Range("H100").Value = 1621.7222
$ws.Range("I100").Value = 1756.5714
$ws.Range("J100").Value = 1535.909
$ws.Range("K100").Value = 1756.5714
$ws.Range("L100").Value = 1535.909
$ws.Range("M100").Value = -1215.5714
$ws.Range("N100").Value = -2617.909
$ws.Range("H116").Value = 4407.5
$ws.Range("I116").Value = 4680
$ws.Range("K116").Value = 4680
$ws.Range("M116").Value = -1238
$ws.Range("H132").Value = 287656.6
$ws.Range("I132").Value = 1989.2858
$ws.Range("K132").Value = 5967.857400000001
$ws.Range("M132").Value = -3437.857400000001
$ws.Range("H137").Value = 852.1539
$ws.Range("I137").Value = 782.0909
$ws.Range("J137").Value = 1237.5
$ws.Range("K137").Value = 2346.2727
$ws.Range("L137").Value = 3712.5
$ws.Range("M137").Value = 203.7273
$ws.Range("N137").Value = -8812.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4703.5796
$ws.Range("I32").Value = 3884.2983
$ws.Range("J32").Value = 8595.166999999999
$ws.Range("K32").Value = 3884.2983
$ws.Range("L32").Value = 8595.166999999999
$ws.Range("M32").Value = -3597.2983
$ws.Range("N32").Value = -9169.166999999999
$ws.Range("H61").Value = 1640.4166
$ws.Range("I61").Value = 1618.6
$ws.Range("K61").Value = 1618.6
$ws.Range("M61").Value = -1406.6
$ws.Range("H88").Value = 2551
$ws.Range("I88").Value = 2232
$ws.Range("J88").Value = 2683.9167
$ws.Range("K88").Value = 2232
$ws.Range("L88").Value = 2683.9167
$ws.Range("M88").Value = -1826
$ws.Range("N88").Value = -3495.9167
$ws.Range("H91").Value = 2551
$ws.Range("I91").Value = 2232
$ws.Range("J91").Value = 2683.9167
$ws.Range("K91").Value = 2232
$ws.Range("L91").Value = 2683.9167
$ws.Range("M91").Value = -828
$ws.Range("N91").Value = -5491.9167
$ws.Range("H102").Value = 3726
$ws.Range("I102").Value = 3544.5454
$ws.Range("J102").Value = 4225
$ws.Range("K102").Value = 3544.5454
$ws.Range("L102").Value = 4225
$ws.Range("M102").Value = -1922.5454
$ws.Range("N102").Value = -7469
$ws.Range("H136").Value = 1640.4166
$ws.Range("I136").Value = 1618.6
$ws.Range("K136").Value = 4855.799999999999
$ws.Range("M136").Value = -2305.799999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H38").Value = 13000
$ws.Range("J38").Value = 13000
$ws.Range("L38").Value = 13000
$ws.Range("N38").Value = -13832
$ws.Range("H94").Value = 2225
$ws.Range("I94").Value = 900
$ws.Range("J94").Value = 2666.6667
$ws.Range("K94").Value = 900
$ws.Range("L94").Value = 2666.6667
$ws.Range("M94").Value = -449
$ws.Range("N94").Value = -3568.6667
$ws.Range("H99").Value = 2901.3125
$ws.Range("I99").Value = 2631
$ws.Range("J99").Value = 3351.8333
$ws.Range("K99").Value = 2631
$ws.Range("L99").Value = 3351.8333
$ws.Range("M99").Value = -1133
$ws.Range("N99").Value = -6347.8333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 2000
$ws.Range("I105").Value = 0
$ws.Range("J105").Value = 2000
$ws.Range("K105").Value = 0
$ws.Range("L105").Value = 2000
$ws.Range("M105").ClearContents()
$ws.Range("N105").Value = -5494

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1486.2354
$ws.Range("I5").Value = 804.4
$ws.Range("J5").Value = 6600
$ws.Range("K5").Value = 2413.2
$ws.Range("L5").Value = 19800
$ws.Range("M5").Value = -2301.2
$ws.Range("N5").Value = -20024
$ws.Range("H113").Value = 605.2727
$ws.Range("I113").Value = 531.6667
$ws.Range("J113").Value = 621.62964
$ws.Range("K113").Value = 1595.0001
$ws.Range("L113").Value = 1864.88892
$ws.Range("M113").Value = 574.9999
$ws.Range("N113").Value = -6204.888919999999
$ws.Range("H122").Value = 1000650.4
$ws.Range("J122").Value = 2500876
$ws.Range("L122").Value = 22507884
$ws.Range("N122").Value = -22512784
$ws.Range("H131").Value = 18377.967
$ws.Range("J131").Value = 1747.1372
$ws.Range("L131").Value = 5241.411599999999
$ws.Range("N131").Value = -15321.4116
$ws.Range("H135").Value = 1486.2354
$ws.Range("I135").Value = 804.4
$ws.Range("J135").Value = 6600
$ws.Range("K135").Value = 7239.599999999999
$ws.Range("L135").Value = 59400
$ws.Range("M135").Value = -4704.599999999999
$ws.Range("N135").Value = -64470

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 8975
$ws.Range("I43").Value = 7333.3335
$ws.Range("J43").Value = 13900
$ws.Range("K43").Value = 7333.3335
$ws.Range("L43").Value = 13900
$ws.Range("M43").Value = -7182.3335
$ws.Range("N43").Value = -14202
$ws.Range("H70").Value = 4351.5
$ws.Range("I70").Value = 4071.261
$ws.Range("K70").Value = 4071.261
$ws.Range("M70").Value = -3801.261
$ws.Range("H73").Value = 4351.5
$ws.Range("I73").Value = 4071.261
$ws.Range("K73").Value = 4071.261
$ws.Range("M73").Value = -3135.261
$ws.Range("H132").Value = 2153
$ws.Range("I132").Value = 1636.1428
$ws.Range("J132").Value = 3359
$ws.Range("K132").Value = 4908.428400000001
$ws.Range("L132").Value = 10077
$ws.Range("M132").Value = -2378.428400000001
$ws.Range("N132").Value = -15137
$ws.Range("H140").Value = 76150.25
$ws.Range("J140").Value = 76150.25
$ws.Range("L140").Value = 76150.25
$ws.Range("N140").Value = -86510.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 2638.9688
$ws.Range("I136").Value = 1461.9584
$ws.Range("J136").Value = 6170
$ws.Range("K136").Value = 4385.8752
$ws.Range("L136").Value = 18510
$ws.Range("M136").Value = -1835.8752
$ws.Range("N136").Value = -23610

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 2247.75
$ws.Range("I100").Value = 330.33334
$ws.Range("K100").Value = 660.66668
$ws.Range("M100").Value = -119.66668
$ws.Range("H136").Value = 1778.8776
$ws.Range("I136").Value = 1671.9117
$ws.Range("J136").Value = 2021.3334
$ws.Range("K136").Value = 5015.7351
$ws.Range("L136").Value = 6064.0002
$ws.Range("M136").Value = -2465.7351
$ws.Range("N136").Value = -11164.0002
